# TagListenersTemplate.xlsx update:
#  - Add a new "before" worksheet as the last tab, containing the
#    "Number of Employees:" label plus two jt:for placeholder cells
#    (mirrors the "before" snapshot added alongside the SheetListener /
#    "before" callback feature work described in the commit message).
#  - Make the new sheet the active tab (was "byClass").

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet (implInstance), so it
# becomes the 5th / final tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "before"

# Write B2 before B1 so the two jt:for strings land in the shared-string
# table in the same order as the source workbook (the "size" expression
# first, then the literal "Will be replaced!" placeholder).
$newSheet.Range("A1").Value = "Number of Employees:"
$newSheet.Range("B2").Value = '<jt:for var="x" start="1" end="1" onProcessed="${boldTagListener}">The above will be replaced by ${employees.size()}</jt:for>'
$newSheet.Range("B1").Value = '<jt:for var="x" start="1" end="1" onProcessed="${boldTagListener}">Will be replaced!</jt:for>'

# Match the original template's auto-sized column A width as closely as
# this host's width quantization allows (source file: width 21.7109375).
$newSheet.Columns("A:A").ColumnWidth = 20.86

# Select/activate the new sheet so it becomes the workbook's active tab
# (tabSelected moves off "byClass" onto "before", activeTab -> 4).
$newSheet.Select()
$newSheet.Activate()
$newSheet.Range("A1").Select()
